$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 594666386

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = "xlPortrait"

$ws.Range("A2").Select()
